# Updates cryptos list values per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.222.99"
$ws.Range("E2").Value = "'  -0.44%  "
$ws.Range("D3").Value = "'3.202.93"
$ws.Range("E3").Value = "'  +0.18%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'607.19"
$ws.Range("E5").Value = "'  +1.61%  "
$ws.Range("D6").Value = "'156.10"
$ws.Range("E6").Value = "'  +0.36%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'3.203.11"
$ws.Range("E8").Value = "'  +0.27%  "
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "'  -2.04%  "
$ws.Range("E10").Value = "'  -0.45%  "
$ws.Range("D11").Value = "'5.65"
$ws.Range("E11").Value = "'  -3.81%  "
$ws.Range("E12").Value = "'  -3.51%  "
$ws.Range("E13").Value = "'  -0.90%  "
$ws.Range("D14").Value = "'38.41"
$ws.Range("E14").Value = "'  -2.36%  "
$ws.Range("D15").Value = "'3.731.36"
$ws.Range("E15").Value = "'  +0.31%  "
$ws.Range("D16").Value = "'66.389.90"
$ws.Range("E16").Value = "'  -0.15%  "
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.207.14"
$ws.Range("E17").Value = "'  +0.24%  "
$ws.Range("B18").Value = "'Polkadot"
$ws.Range("C18").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.24"
$ws.Range("E18").Value = "'  -3.72%  "
$ws.Range("E19").Value = "'  +1.21%  "
$ws.Range("D20").Value = "'505.39"
$ws.Range("E21").Value = "'  -1.08%  "
$ws.Range("E23").Value = "'  -2.06%  "
$ws.Range("D24").Value = "'14.62"
$ws.Range("E24").Value = "'  -2.34%  "
$ws.Range("D25").Value = "'85.09"
$ws.Range("E25").Value = "'  -1.00%  "
$ws.Range("E26").Value = "'  -0.11%  "
$ws.Range("D27").Value = "'2.99"
$ws.Range("E27").Value = "'  -0.54%  "
$ws.Range("D28").Value = "'9.03"
$ws.Range("E28").Value = "'  -2.74%  "
$ws.Range("E29").Value = "'  -1.08%  "
$ws.Range("D30").Value = "'0.131"
$ws.Range("E30").Value = "'  +44.27%  "
$ws.Range("D31").Value = "'2.91"
$ws.Range("E31").Value = "'  -1.01%  "
$ws.Range("E32").Value = "'  -1.97%  "
$ws.Range("D33").Value = "'28.17"
$ws.Range("E34").Value = "'  +0.07%  "
$ws.Range("E35").Value = "'  -5.05%  "
$ws.Range("E36").Value = "'  -2.13%  "
$ws.Range("D37").Value = "'55.34"
$ws.Range("E37").Value = "'  +0.74%  "
$ws.Range("D38").Value = "'499.67"
$ws.Range("E38").Value = "'  -2.25%  "
$ws.Range("E39").Value = "'  +12.23%  "
$ws.Range("E40").Value = "'  +2.24%  "
$ws.Range("E41").Value = "'  +5.12%  "
$ws.Range("E42").Value = "'  -1.86%  "
$ws.Range("E43").Value = "'  -2.53%  "
$ws.Range("D44").Value = "'0.295"
$ws.Range("E44").Value = "'  -2.02%  "
$ws.Range("D45").Value = "'2.909.98"
$ws.Range("E45").Value = "'  +0.27%  "
$ws.Range("E46").Value = "'  -1.17%  "
$ws.Range("D47").Value = "'28.14"
$ws.Range("E47").Value = "'  -1.54%  "
$ws.Range("E48").Value = "'  +1.50%  "
$ws.Range("E49").Value = "'  -0.05%  "
$ws.Range("E50").Value = "'  -1.14%  "
$ws.Range("D51").Value = "'122.13"
$ws.Range("E51").Value = "'  +0.39%  "
